$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.855.19'
$ws.Range('E2').Value = '  -1.23%  '
$ws.Range('D3').Value = '3.133.61'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.65'
$ws.Range('D5').Style = $ws.Range('A1').Style
$ws.Range('E5').Value = '  -2.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.89'
$ws.Range('D6').Style = $ws.Range('A1').Style
$ws.Range('E6').Value = '  -3.78%  '
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('D8').Value = '3.129.06'
$ws.Range('E8').Value = '  -0.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.525'
$ws.Range('D9').Style = $ws.Range('A1').Style
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').Value = '  -2.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.33'
$ws.Range('D11').Style = $ws.Range('A1').Style
$ws.Range('E11').Value = '  -1.90%  '
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000253'
$ws.Range('D13').Style = $ws.Range('A1').Style
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.53'
$ws.Range('D14').Style = $ws.Range('A1').Style
$ws.Range('E14').Value = '  -2.92%  '
$ws.Range('D15').Value = '3.650.68'
$ws.Range('E16').Value = '  +2.58%  '
$ws.Range('D17').Value = '63.764.00'
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('D18').Value = '3.135.98'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.77'
$ws.Range('D19').Style = $ws.Range('A1').Style
$ws.Range('E19').Value = '  -1.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '481.50'
$ws.Range('D20').Style = $ws.Range('A1').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.52'
$ws.Range('D21').Style = $ws.Range('A1').Style
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.705'
$ws.Range('D22').Style = $ws.Range('A1').Style
$ws.Range('E22').Value = '  -2.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.68'
$ws.Range('D23').Style = $ws.Range('A1').Style
$ws.Range('E23').Value = '  -3.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.62'
$ws.Range('D24').Style = $ws.Range('A1').Style
$ws.Range('E24').Value = '  +4.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.11'
$ws.Range('D25').Style = $ws.Range('A1').Style
$ws.Range('E25').Value = '  -4.88%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -2.50%  '
$ws.Range('E28').Value = '  -6.08%  '
$ws.Range('E29').Value = '  -2.49%  '
$ws.Range('E30').Value = '  -2.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.16'
$ws.Range('D31').Style = $ws.Range('A1').Style
$ws.Range('E31').Value = '  +2.47%  '
$ws.Range('E32').Value = '  -7.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.59'
$ws.Range('D34').Style = $ws.Range('A1').Style
$ws.Range('E34').Value = '  -3.47%  '
$ws.Range('E35').Value = '  -2.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.00'
$ws.Range('D36').Style = $ws.Range('A1').Style
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.54'
$ws.Range('D37').Style = $ws.Range('A1').Style
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('D38').Value = '0.0₃0735'
$ws.Range('E38').Value = '  -6.37%  '
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '429.80'
$ws.Range('D40').Style = $ws.Range('A1').Style
$ws.Range('E40').Value = '  -6.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.84'
$ws.Range('D41').Style = $ws.Range('A1').Style
$ws.Range('E41').Value = '  -10.62%  '
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.30'
$ws.Range('D43').Style = $ws.Range('A1').Style
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').Value = '2.876.80'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.259'
$ws.Range('D45').Style = $ws.Range('A1').Style
$ws.Range('E45').Value = '  -3.11%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.16'
$ws.Range('D46').Style = $ws.Range('A1').Style
$ws.Range('E46').Value = '  -6.68%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.37'
$ws.Range('D47').Style = $ws.Range('A1').Style
$ws.Range('E47').Value = '  -3.20%  '
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.52'
$ws.Range('D50').Style = $ws.Range('A1').Style
$ws.Range('E50').Value = '  -3.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.45'
$ws.Range('D51').Style = $ws.Range('A1').Style
$ws.Range('E51').Value = '  +0.43%  '
